# The commit only touches the library/architecture description text
# body on slide 3 (placeholder shape "Объект 3"). Locate it by content
# instead of a hard-coded index so the script is resilient to shape
# ordering.
$p = $ppt.ActivePresentation
$tr = $null
for ($si = 1; $si -le $p.Slides.Count -and $tr -eq $null; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text.IndexOf("werkzeug.security") -ge 0) {
                $tr = $shape.TextFrame.TextRange
                break
            }
        }
    }
}

# --- Edit 1 -----------------------------------------------------------
# ", game, werkzeug.security(" -> ", erkzeug.security("
#   a) drop the word "game" (keep surrounding ", " / ", ")
#   b) drop the leading "w" of "werkzeug.security"
$full = $tr.Text
$idx = $full.IndexOf(", game, ")
$start = $idx + 1
$len = ", game, ".Length
$tr.Characters($start, $len).Text = ", "

$full = $tr.Text
$idx = $full.IndexOf("werkzeug.security")
$start = $idx + 1
$tr.Characters($start, 1).Text = ""

# --- Edit 2 -----------------------------------------------------------
# "...рендерит HTML из templates.Статика (CSS/JS)..."
#   -> "...рендерит HTML. Статика (CSS/JS)..."
# i.e. "HTML " -> "HTML.", then drop "из " and "templates.", keeping a
# freshly-split " Статика " run followed by its own "(" run.
$full = $tr.Text
$idx = $full.IndexOf("HTML ")
$start = $idx + 1
$len = "HTML ".Length
$tr.Characters($start, $len).Text = "HTML."

$full = $tr.Text
$idx = $full.IndexOf("из templates.Статика (")
$start = $idx + 1
$len = "из templates.Статика (".Length
$tr.Characters($start, $len).Text = " Статика ("

# Re-stamp the trailing "(" on its own so it becomes a distinct run
# (matching the "(" run left behind once "Статика (" was split off in
# the source edit) instead of staying merged with " Статика ".
$full = $tr.Text
$idx = $full.IndexOf(" Статика (")
$start = $idx + 1 + " Статика ".Length
$tr.Characters($start, 1).Text = "("
